$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the same date value (45177) for every data
# row (2 through 224). The update bumps that date by one day (45178).
$ws.Range("C2:C224").Value = 45178
